# "Generate Report for handback"
#
# The 56bec4ee-7eef-4f4f-8950-234bdbffa32a.md file has now been handed
# back (it was previously only "Ready for handoff"). Update the status
# on the Overview sheet and on each locale sheet, and stamp the new
# "Latest Handback DateTime" for each locale.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the 56bec4ee... file (zh-cn / de-de status) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: row 3 status + new handback datetime ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $statusHandedBack
$zhcn.Range("G3").Value = "2016-01-22 02:31:40"

# --- de-de sheet: row 3 status + new handback datetime ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $statusHandedBack
$dede.Range("G3").Value = "2016-01-22 02:32:03"
